# Replace the FY2020 file/table mapping list with the refreshed
# FY2021_Q2 file names (pulled from the new ArcGIS/SharePoint account),
# while keeping the existing table_name values (pw_disclosure is now
# reused for two source files).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("file_name", "table_name"),
    @("H-2A_Addendum_A_Disclosure_Data_FY2021_Q2.xls", "h2a_addendum_a"),
    @("H-2B_Appendix_A_FY2021_Q2.xls",                 "h2b_appendix_a"),
    @("H-2B_Appendix_C_FY2021_Q2.xls",                 "h2b_appendix_c"),
    @("H-2B_Appendix_D_FY2021_Q2.xls",                 "h2b_appendix_d"),
    @("PW_Disclosure_Data_FY2021_Q2.xls",               "pw_disclosure"),
    @("PW_Disclosure_Data_FY2021_Q2_2.xls",             "pw_disclosure"),
    @("PW_Worksites_FY2021_Q2.xls",                     "pw_worksites")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Clear anything that might linger below the new table (defensive;
# the table is the same size as before so this is a no-op today).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -gt $data.Length) {
    $ws.Range($ws.Cells.Item($data.Length + 1, 1), $ws.Cells.Item($lastRow, 2)).ClearContents()
}

# Mirror the new selection left behind in the workbook (A2:B8, active cell A2)
$ws.Range("A2:B8").Select()

$wb.Save()
